$wb = $excel.ActiveWorkbook

# --- "model" sheet: insert a new field row (REGID / test / FALSE) before
# the existing REGIDC row, shifting REGIDC/CONT/comsup/parma/moma down one
# row each (fixes the linked table lookup field used by the query below).
$model = $wb.Worksheets.Item("model")
$model.Rows.Item(24).Insert()
$model.Cells.Item(24, 1).Value = "REGID"
$model.Cells.Item(24, 2).Value = "test"
$model.Cells.Item(24, 3).Value = $false

# --- "queries" sheet: point the linked_visitdate query's selection /
# selectionArgs at the new REGID field instead of REGIDC.
$queries = $wb.Worksheets.Item("queries")
$queries.Range("E4").Value = "REGID = ?"
$queries.Range("F4").Value = "[data('REGID')]"

# --- restore view state (selection/scroll position) on the sheets that
# were touched while editing, then leave "queries" as the active tab.
$survey = $wb.Worksheets.Item("survey")
$survey.Activate()
$survey.Range("D9").Select()

$choices = $wb.Worksheets.Item("choices")
$choices.Activate()
$choices.Range("D13").Select()

$model.Activate()
$model.Range("B25").Select()

$queries.Activate()
$queries.Range("F4").Select()
